# Adapt column header formatting to respective input file names (#7)
#
# 1. Rename the "_old" / "_new" header-name suffixes (columns A1:J1 and
#    L1:U1) to "_FV2210" / "_FV2304" respectively (column K1 "diff" is
#    left untouched).
# 2. Turn the header row + data range into a native Excel Table
#    ("Table1") spanning A1:U64.
# 3. Freeze the header row (split after row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename header cells
# ---------------------------------------------------------------------
$headersFV2210 = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

$headersFV2304 = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headersFV2210.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2210[$i]
}
# Column 11 (K1) stays "diff" - no change needed.
for ($i = 0; $i -lt $headersFV2304.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2304[$i]
}

# ---------------------------------------------------------------------
# 2) Convert A1:U64 into a Table, while preserving the header row's
#    existing formatting exactly (bold / shaded / bordered / centered /
#    wrapped) instead of letting Excel overwrite it with a table-style
#    dxf override. We stash a copy of the header formatting on an unused
#    scratch row, let Excel build the table (which strips/replaces the
#    header formatting), paste the formatting back, then remove the
#    scratch row again.
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRow = 100
$scratch = $ws.Range("A" + $scratchRow + ":U" + $scratchRow)

$headerRange.Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null    # xlPasteFormats

$headerRange.ClearFormats() | Out-Null

$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U64"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

$scratch.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Rows.Item($scratchRow).Delete() | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Freeze the header row (row 1).
# ---------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
